$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New data rows to append (A=DATA serial, B..J values)
$data = @(
    @(44377, 12872, 117, 5970, 18959, 5620, 214, 19, 195, 136),
    @(44378, 12954, 98, 5993, 19045, 5651, 206, 19, 187, 136),
    @(44379, 12987, 111, 6004, 19102, 5677, 189, 18, 171, 138),
    @(44380, 13046, 84, 6019, 19149, 5702, 179, 18, 161, 138),
    @(44381, 13088, 53, 6025, 19166, 5718, 169, 17, 152, 138)
)

$startRow = 364
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 10; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$col - 1]
    }
}

# Reuse the existing formatting (style indices) from row 2 by copying formats only,
# so no new cellXfs entries are created.
$ws.Range("A2").Copy()
$ws.Range("A364:A368").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B2:J2").Copy()
$ws.Range("B364:J368").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Update selection to match the post-edit state
$ws.Range("K371").Select()
